$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update kategori_id values in column A (rows 2-6)
$ws.Range("A2").Value = 3
$ws.Range("A3").Value = 6
$ws.Range("A4").Value = 5
$ws.Range("A5").Value = 5
$ws.Range("A6").Value = 4

# Update the active selection/cell to C15
$ws.Range("C15").Select()
